$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the error message and date on row 4 (F4 = error column, G4 = date column)
$ws.Range("F4").Value = "Coould not fetch interest Rates. Invalid date 2022-09-05, valid date format exemple: 05/09/2022"
$ws.Range("G4").Value = "2022-09-05 21:18:10"
